# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.862.25"
$ws.Range("E2").Value = "  +4.48%  "
$ws.Range("D3").Value = "3.340.57"
$ws.Range("E3").Value = "  +4.48%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'556.05"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").Value = "'152.22"
$ws.Range("E6").Value = "  +4.80%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  +3.96%  "
$ws.Range("D11").Value = "'0.438"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").Value = "3.910.44"
$ws.Range("E12").Value = "  +4.34%  "
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "'0.0000182"
$ws.Range("E14").Value = "  +3.92%  "
$ws.Range("D15").Value = "'26.89"
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").Value = "62.855.61"
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.320.25"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'6.48"
$ws.Range("E18").Value = "  +4.55%  "
$ws.Range("D19").Value = "'13.75"
$ws.Range("E19").Value = "  +4.91%  "
$ws.Range("D20").Value = "'8.45"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "'388.32"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D24").Value = "'70.63"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").Value = "'0.180"
$ws.Range("E25").Value = "  +4.95%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "0.0₃0965"
$ws.Range("E27").Value = "  +6.86%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("D30").Value = "'6.39"
$ws.Range("E30").Value = "  +3.56%  "
$ws.Range("D31").Value = "'22.99"
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("E32").Value = "  +2.66%  "
$ws.Range("D33").Value = "'1.30"
$ws.Range("E33").Value = "  +6.35%  "
$ws.Range("D34").Value = "'6.69"
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.49"
$ws.Range("E35").Value = "  +9.95%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'160.55"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").Value = "'1.89"
$ws.Range("E37").Value = "  +11.09%  "
$ws.Range("D38").Value = "'27.10"
$ws.Range("E38").Value = "  +5.12%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.0738"
$ws.Range("E39").Value = "  +3.23%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.819.66"
$ws.Range("E40").Value = "  +1.68%  "
$ws.Range("D41").Value = "'0.0312"
$ws.Range("E41").Value = "  +8.45%  "
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'40.65"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.746"
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("D46").Value = "3.381.47"
$ws.Range("E46").Value = "  +4.38%  "
$ws.Range("D47").Value = "'21.89"
$ws.Range("E47").Value = "  +6.43%  "
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("D50").Value = "'0.800"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "'282.24"
$ws.Range("E51").Value = "  +6.70%  "
